# Insert two new data rows at the top of the "Vega Monumental Concepción - Limón"
# date block (rows 276-277), pushing the existing rows 276-340 down to 278-342.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 276.
$ws.Range("A276:A277").EntireRow.Insert()

# --- New row 276 ---
$ws.Range("A276").Value = 11
$ws.Range("B276").Value = "Vega Monumental Concepción"
$ws.Range("C276").Value = "Bíobío"
$ws.Range("D276").Value = 44551
$ws.Range("E276").Value = 8
$ws.Range("F276").Value = "Fruta"
$ws.Range("G276").Value = 100102
$ws.Range("H276").Value = "Cítricos"
$ws.Range("I276").Value = 100102003
$ws.Range("J276").Value = "Limón"
$ws.Range("K276").Value = "Sin especificar"
$ws.Range("L276").Value = "1a amarillo"
$ws.Range("M276").Value = 350
$ws.Range("N276").Value = 15000
$ws.Range("O276").Value = 16000
$ws.Range("P276").Value = 15571
$ws.Range("Q276").Value = "$/malla 18 kilos"
$ws.Range("R276").Value = "Región Metropolitana"
$ws.Range("S276").Value = 865
$ws.Range("T276").Value = 18

# --- New row 277 ---
$ws.Range("A277").Value = 11
$ws.Range("B277").Value = "Vega Monumental Concepción"
$ws.Range("C277").Value = "Bíobío"
$ws.Range("D277").Value = 44551
$ws.Range("E277").Value = 8
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100102
$ws.Range("H277").Value = "Cítricos"
$ws.Range("I277").Value = 100102003
$ws.Range("J277").Value = "Limón"
$ws.Range("K277").Value = "Sin especificar"
$ws.Range("L277").Value = "1a plateado"
$ws.Range("M277").Value = 220
$ws.Range("N277").Value = 14000
$ws.Range("O277").Value = 15000
$ws.Range("P277").Value = 14455
$ws.Range("Q277").Value = "$/malla 18 kilos"
$ws.Range("R277").Value = "Provincia de Melipilla"
$ws.Range("S277").Value = 803
$ws.Range("T277").Value = 18
